# Progress update as of 04-Nov-2025:
#  - "PERIOD TO EXPIRE" (col H) drops by 1 day for every training row.
#  - "LAST UPDATE" (col I) moves from 03-Nov-2025 to 04-Nov-2025.
# Applies to rows 3-12 of the "Training Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$rows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12)

# Scratch cell well outside the used range (A1:K12) used to coerce the
# TEXT() formula result into a genuine literal string value via
# Copy/PasteSpecial (values only). A plain `.Value = "04-Nov-2025"`
# assignment gets auto-detected by Excel as a date and silently reformats
# the cell (new number format / style) -- going through TEXT() + paste
# values keeps the cell a plain string and leaves its style untouched,
# matching how the sheet already stored these dates (plain text).
$scratch = $ws.Range("Z1")
$scratch.Formula = "=TEXT(""04-Nov-2025"",""@"")"
$scratch.Copy()

foreach ($r in $rows) {
    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value = $hCell.Value() - 1

    $iCell = $ws.Cells.Item($r, 9)
    $iCell.PasteSpecial(-4163)
}

$scratch.Clear()
$excel.CutCopyMode = $false
